$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted before the current row 84,
# pushing the existing rows 84:109 down to 85:110 (dimension grows to R110).
$ws.Rows.Item(84).Insert()

# Populate the newly inserted row 84 with the new record's values.
$ws.Range("A84").Value = 4
$ws.Range("B84").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C84").Value = "Los Lagos"
$ws.Range("D84").Value = 44876
$ws.Range("E84").Value = 10
$ws.Range("F84").Value = 100112031
$ws.Range("G84").Value = "Poroto verde"
$ws.Range("H84").Value = "Magnum"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 35
$ws.Range("K84").Value = 40000
$ws.Range("L84").Value = 40000
$ws.Range("M84").Value = 40000
$ws.Range("N84").Value = "$/malla 25 kilos"
$ws.Range("O84").Value = "Perú"
$ws.Range("P84").Value = 1600
$ws.Range("Q84").Value = 25
$ws.Range("R84").Value = "Hortaliza"
